$wb = $excel.ActiveWorkbook

# Metadata sheet: row 4 = Name property -> set value to "CapacitesavoirfaireVs"
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B4").Value = "CapacitesavoirfaireVs"

# Update the Date value (row 8, column B) to the new timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
